# Weekly update: insert 5 new price rows (one new reporting week) for
# "Comercializadora del Agro de Limarí" / Mandarina, pushing the existing
# historical rows (385:420) down to (390:425).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows starting at row 385 (shifts old 385:420 -> 390:425)
$ws.Range("A385:A389").EntireRow.Insert()

# Fixed/common columns for every row in this market+product block
$mercadoId  = 2
$mercado    = "Comercializadora del Agro de Limarí"
$region     = "Coquimbo"
$codreg     = 4
$tipo       = "Fruta"
$productoId = 100102
$producto   = "Cítricos"
$categoriaId = 100102004
$categoria  = "Mandarina"
$origen     = "Provincia de Limarí"

# New rows: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Precio$/Kg, Kg/unidad
$newRows = @(
    @(385, 44769, "Clemenuless", "Especial", 400, 6500, 7000, 6750, "`$/bandeja 10 kilos", 675, 10),
    @(386, 44769, "Clemenuless", "Primera",  400, 5500, 6000, 5750, "`$/bandeja 10 kilos", 575, 10),
    @(387, 44769, "Clemenuless", "Primera",  20,  155000, 160000, 157500, "`$/bins (450 kilos)", 350, 450),
    @(388, 44769, "Clemenuless", "Segunda",  400, 4500, 5000, 4750, "`$/bandeja 10 kilos", 475, 10),
    @(389, 44769, "Clemenuless", "Segunda",  20,  115000, 120000, 117500, "`$/bins (450 kilos)", 261, 450)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r[1]
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r[2]
    $ws.Cells.Item($row, 12).Value = $r[3]
    $ws.Cells.Item($row, 13).Value = $r[4]
    $ws.Cells.Item($row, 14).Value = $r[5]
    $ws.Cells.Item($row, 15).Value = $r[6]
    $ws.Cells.Item($row, 16).Value = $r[7]
    $ws.Cells.Item($row, 17).Value = $r[8]
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r[9]
    $ws.Cells.Item($row, 20).Value = $r[10]
}
